# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated output (per commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Map: sheet name -> { row -> new value }
$updates = @{
    "展览" = @{
        3  = 3973
        4  = 2333
        5  = 464
        11 = 47
        13 = 1475
        15 = 2732
    }
    "全部类型" = @{
        3  = 3973
        4  = 2333
        5  = 464
        12 = 47
        16 = 1475
        18 = 2732
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
